$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths
$ws.Columns.Item(1).ColumnWidth = 16.0
$ws.Columns.Item(2).ColumnWidth = 29.5
$ws.Columns.Item(3).ColumnWidth = 28.0
$ws.Columns.Item(4).ColumnWidth = 34.333333333333336
$ws.Columns.Item(5).ColumnWidth = 28.166666666666668
$ws.Columns.Item(6).ColumnWidth = 26.833333333333332
$ws.Columns.Item(7).ColumnWidth = 33.166666666666664
$ws.Columns.Item(8).ColumnWidth = 28.333333333333332
$ws.Columns.Item(9).ColumnWidth = 27.166666666666668

# Cell values
$ws.Range("B2").Value = 34.33393724753568
$ws.Range("C2").Value = 15.534752777185458
$ws.Range("D2").Value = 0.45246056882976526
$ws.Range("E2").Value = 29.84715836273299
$ws.Range("F2").Value = 14.584827266122947
$ws.Range("G2").Value = 0.4886504466815001
$ws.Range("H2").Value = 338
$ws.Range("I2").Value = 300
$ws.Range("B3").Value = 40.264225817871036
$ws.Range("C3").Value = 21.184180140110044
$ws.Range("D3").Value = 0.5261290813322325
$ws.Range("E3").Value = 35.85611015272576
$ws.Range("F3").Value = 19.80771824360334
$ws.Range("G3").Value = 0.5524223949344814
$ws.Range("H3").Value = 342
$ws.Range("I3").Value = 309
$ws.Range("B4").Value = 45.675747208180894
$ws.Range("C4").Value = 26.324666404641626
$ws.Range("D4").Value = 0.5763379476784272
$ws.Range("E4").Value = 42.04334371929386
$ws.Range("F4").Value = 25.177068083047473
$ws.Range("G4").Value = 0.5988360072201778
$ws.Range("H4").Value = 346.5
$ws.Range("I4").Value = 317
$ws.Range("B5").Value = 51.49264138024465
$ws.Range("C5").Value = 31.678939398108426
$ws.Range("D5").Value = 0.6152129420624783
$ws.Range("E5").Value = 47.83648288471331
$ws.Range("F5").Value = 30.42367067958894
$ws.Range("G5").Value = 0.6359930505950965
$ws.Range("H5").Value = 351
$ws.Range("I5").Value = 323.5
$ws.Range("B6").Value = 57.89644856269853
$ws.Range("C6").Value = 38.67046252473081
$ws.Range("D6").Value = 0.6679246047856099
$ws.Range("E6").Value = 53.50851768750377
$ws.Range("F6").Value = 35.80562983180617
$ws.Range("G6").Value = 0.6691575730225866
$ws.Range("H6").Value = 357.5
$ws.Range("I6").Value = 329
$ws.Range("B7").Value = 63.64781193965802
$ws.Range("C7").Value = 43.378568863484894
$ws.Range("D7").Value = 0.6815406145400631
$ws.Range("E7").Value = 59.48161734192175
$ws.Range("F7").Value = 41.58269405078744
$ws.Range("G7").Value = 0.6990847913861378
$ws.Range("H7").Value = 362
$ws.Range("I7").Value = 335.5
$ws.Range("B8").Value = 69.59054454925227
$ws.Range("C8").Value = 49.26845483111522
$ws.Range("D8").Value = 0.7079762796833093
$ws.Range("E8").Value = 65.78784392096837
$ws.Range("F8").Value = 47.644760945303815
$ws.Range("G8").Value = 0.7242183069951338
$ws.Range("H8").Value = 367.5
$ws.Range("I8").Value = 341.5
$ws.Range("B9").Value = 75.79133736321782
$ws.Range("C9").Value = 55.34887193389558
$ws.Range("D9").Value = 0.7302796580649448
$ws.Range("E9").Value = 71.7318861922717
$ws.Range("F9").Value = 53.2949507769771
$ws.Range("G9").Value = 0.7429743396698668
$ws.Range("H9").Value = 372.5
$ws.Range("I9").Value = 347
$ws.Range("B10").Value = 81.6612686502326
$ws.Range("C10").Value = 61.01429026060374
$ws.Range("D10").Value = 0.7471631444024834
$ws.Range("E10").Value = 77.31706821749636
$ws.Range("F10").Value = 58.66959184723146
$ws.Range("G10").Value = 0.7588181135139693
$ws.Range("H10").Value = 377
$ws.Range("I10").Value = 352
$ws.Range("B11").Value = 87.76505904031158
$ws.Range("C11").Value = 66.94212917516893
$ws.Range("D11").Value = 0.7627423704508827
$ws.Range("E11").Value = 83.15411552205771
$ws.Range("F11").Value = 64.3889649405951
$ws.Range("G11").Value = 0.7743328701934794
$ws.Range("H11").Value = 381.5
$ws.Range("I11").Value = 357
